$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Lrpap1"
$ws.Range("C2").Value = "Vldlr"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 5.652167666666667
$ws.Range("H2").Value = 16.956503
$ws.Range("I2").Value = 0.1860329065948871
$ws.Range("J2").Value = 0.1860329065948871
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.06089466666666667
$ws.Range("N2").Value = 0.182684
$ws.Range("O2").Value = 0.001903591634475228
$ws.Range("P2").Value = 0.001903591634475228
$ws.Range("Q2").Value = 0.3441868660057779
$ws.Range("R2").Value = 3.097681794052
$ws.Range("S2").Value = 0.0003541306847311385
$ws.Range("T2").Value = 0.0003541306847311385

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Lrpap1"
$ws.Range("C3").Value = "Vldlr"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 5.652167666666667
$ws.Range("H3").Value = 16.956503
$ws.Range("I3").Value = 0.1860329065948871
$ws.Range("J3").Value = 0.1860329065948871
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 28.046323
$ws.Range("N3").Value = 84.138969
$ws.Range("O3").Value = 0.8767392739472014
$ws.Range("P3").Value = 0.8767392739472013
$ws.Range("Q3").Value = 158.5225200294897
$ws.Range("R3").Value = 1426.702680265407
$ws.Range("S3").Value = 0.1631023554582888
$ws.Range("T3").Value = 0.1631023554582888

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Lrpap1"
$ws.Range("C4").Value = "Vldlr"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 5.652167666666667
$ws.Range("H4").Value = 16.956503
$ws.Range("I4").Value = 0.1860329065948871
$ws.Range("J4").Value = 0.1860329065948871
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.882136333333333
$ws.Range("N4").Value = 11.646409
$ws.Range("O4").Value = 0.1213571344183235
$ws.Range("P4").Value = 0.1213571344183235
$ws.Range("Q4").Value = 21.94248546085856
$ws.Range("R4").Value = 197.482369147727
$ws.Range("S4").Value = 0.02257642045186712
$ws.Range("T4").Value = 0.02257642045186712

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Lrpap1"
$ws.Range("C5").Value = "Vldlr"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 16.407289
$ws.Range("H5").Value = 49.221867
$ws.Range("I5").Value = 0.5400221369958743
$ws.Range("J5").Value = 0.5400221369958743
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.06089466666666667
$ws.Range("N5").Value = 0.182684
$ws.Range("O5").Value = 0.001903591634475228
$ws.Range("P5").Value = 0.001903591634475228
$ws.Range("Q5").Value = 0.999116394558667
$ws.Range("R5").Value = 8.992047551028001
$ws.Range("S5").Value = 0.001027981622416782
$ws.Range("T5").Value = 0.001027981622416782

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Lrpap1"
$ws.Range("C6").Value = "Vldlr"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 16.407289
$ws.Range("H6").Value = 49.221867
$ws.Range("I6").Value = 0.5400221369958743
$ws.Range("J6").Value = 0.5400221369958743
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 28.046323
$ws.Range("N6").Value = 84.138969
$ws.Range("O6").Value = 0.8767392739472014
$ws.Range("P6").Value = 0.8767392739472013
$ws.Range("Q6").Value = 460.1641268483471
$ws.Range("R6").Value = 4141.477141635123
$ws.Range("S6").Value = 0.473458616305179
$ws.Range("T6").Value = 0.4734586163051789

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Lrpap1"
$ws.Range("C7").Value = "Vldlr"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 16.407289
$ws.Range("H7").Value = 49.221867
$ws.Range("I7").Value = 0.5400221369958743
$ws.Range("J7").Value = 0.5400221369958743
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.882136333333333
$ws.Range("N7").Value = 11.646409
$ws.Range("O7").Value = 0.1213571344183235
$ws.Range("P7").Value = 0.1213571344183235
$ws.Range("Q7").Value = 63.69533275840034
$ws.Range("R7").Value = 573.2579948256031
$ws.Range("S7").Value = 0.06553553906827861
$ws.Range("T7").Value = 0.06553553906827861

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Lrpap1"
$ws.Range("C8").Value = "Vldlr"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 8.323166333333333
$ws.Range("H8").Value = 24.969499
$ws.Range("I8").Value = 0.2739449564092387
$ws.Range("J8").Value = 0.2739449564092387
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.06089466666666667
$ws.Range("N8").Value = 0.182684
$ws.Range("O8").Value = 0.001903591634475228
$ws.Range("P8").Value = 0.001903591634475228
$ws.Range("Q8").Value = 0.5068364394795556
$ws.Range("R8").Value = 4.561527955316
$ws.Range("S8").Value = 0.0005214793273273078
$ws.Range("T8").Value = 0.0005214793273273078

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Lrpap1"
$ws.Range("C9").Value = "Vldlr"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 8.323166333333333
$ws.Range("H9").Value = 24.969499
$ws.Range("I9").Value = 0.2739449564092387
$ws.Range("J9").Value = 0.2739449564092387
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 28.046323
$ws.Range("N9").Value = 84.138969
$ws.Range("O9").Value = 0.8767392739472014
$ws.Range("P9").Value = 0.8767392739472013
$ws.Range("Q9").Value = 233.4342113673923
$ws.Range("R9").Value = 2100.907902306531
$ws.Range("S9").Value = 0.2401783021837336
$ws.Range("T9").Value = 0.2401783021837336

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Lrpap1"
$ws.Range("C10").Value = "Vldlr"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 8.323166333333333
$ws.Range("H10").Value = 24.969499
$ws.Range("I10").Value = 0.2739449564092387
$ws.Range("J10").Value = 0.2739449564092387
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.882136333333333
$ws.Range("N10").Value = 11.646409
$ws.Range("O10").Value = 0.1213571344183235
$ws.Range("P10").Value = 0.1213571344183235
$ws.Range("Q10").Value = 32.31166643101011
$ws.Range("R10").Value = 290.804997879091
$ws.Range("S10").Value = 0.03324517489817774
$ws.Range("T10").Value = 0.03324517489817774
